# Actualización desde MV -datos-
# Adds the "01-07-2021" period (row 144) and revises several previously
# reported monthly figures (rows 138-143) in the "Operaciones no
# financieras" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revisions to existing rows 138-143 -----------------------------------

# Row 138
$ws.Range("B138").Value = 4263723
$ws.Range("D138").Value = 46097
$ws.Range("J138").Value = 3583204
$ws.Range("M138").Value = 365940
$ws.Range("Q138").Value = 680519
$ws.Range("W138").Value = 3790720
$ws.Range("X138").Value = 473232

# Row 139
$ws.Range("J139").Value = 3446985
$ws.Range("M139").Value = 30593
$ws.Range("Q139").Value = 110547
$ws.Range("W139").Value = 3882671
$ws.Range("X139").Value = -324827

# Row 140
$ws.Range("B140").Value = 4246316
$ws.Range("D140").Value = 105791
$ws.Range("I140").Value = 387611
$ws.Range("J140").Value = 5258995
$ws.Range("M140").Value = 519236
$ws.Range("Q140").Value = -1012679
$ws.Range("V140").Value = 4247580
$ws.Range("W140").Value = 5903745
$ws.Range("X140").Value = -1656164

# Row 141
$ws.Range("B141").Value = 6816900
$ws.Range("D141").Value = 166597
$ws.Range("J141").Value = 4818992
$ws.Range("M141").Value = 39203
$ws.Range("Q141").Value = 1997908
$ws.Range("V141").Value = 6817230
$ws.Range("W141").Value = 5406209
$ws.Range("X141").Value = 1411020

# Row 142
$ws.Range("J142").Value = 5780773
$ws.Range("M142").Value = 31776
$ws.Range("Q142").Value = -2790265
$ws.Range("W142").Value = 6338450
$ws.Range("X142").Value = -3347299

# Row 143
$ws.Range("B143").Value = 3732846
$ws.Range("I143").Value = 174624
$ws.Range("J143").Value = 5319558
$ws.Range("M143").Value = 22751
$ws.Range("Q143").Value = -1586712
$ws.Range("V143").Value = 3733041
$ws.Range("W143").Value = 5991127
$ws.Range("X143").Value = -2258086

# --- New row 144 (period 01-07-2021) ---------------------------------------

# Column A holds the period label as text ("01-07-2021"). A bare .Value
# assignment gets auto-parsed into a date serial by Excel's type inference,
# so force it in as text via a leading apostrophe, then strip the resulting
# cell style back to the sheet's default (matching the rest of column A,
# which carries no explicit style).
$ws.Range("A144").Value = "'01-07-2021"
$ws.Range("A144").Style = "Normal"

$ws.Range("B144").Value = 4139698
$ws.Range("C144").Value = 3254655
$ws.Range("D144").Value = 291577
$ws.Range("E144").Value = 231852
$ws.Range("F144").Value = 3554
$ws.Range("G144").Value = 37084
$ws.Range("H144").Value = 78197
$ws.Range("I144").Value = 242779
$ws.Range("J144").Value = 7004507
$ws.Range("K144").Value = 876783
$ws.Range("L144").Value = 372509
$ws.Range("M144").Value = 335936
$ws.Range("N144").Value = 4692589
$ws.Range("O144").Value = 703576
$ws.Range("P144").Value = 23115
$ws.Range("Q144").Value = -2864809
$ws.Range("R144").Value = 569551
$ws.Range("S144").Value = 574
$ws.Range("T144").Value = 297190
$ws.Range("U144").Value = 272935
$ws.Range("V144").Value = 4140273
$ws.Range("W144").Value = 7574632
$ws.Range("X144").Value = -3434360
